$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.767.65"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.523.07"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.87"
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.77"
$ws.Range("E6").Value = "  +5.36%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.522.20"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.80"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.978.67"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.640.54"
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.521.81"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.97"
$ws.Range("E19").Value = "  +5.45%  "
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.76"
$ws.Range("E21").Value = "  +5.19%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.64"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.17"
$ws.Range("E26").Value = "  +3.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0991"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "541.95"
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.79"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.81"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.561"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0280"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.18"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0756"
$ws.Range("E51").Value = "  +1.01%  "
